$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create "Challenge 5" by copying the existing "Challenge 4" sheet (it has
#    the closest layout/styles to the new sheet) and placing it right after.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("Challenge 4")
$src.Copy([Type]::Missing, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "Challenge 5"
$ws.Range("F46").Select()

# ---------------------------------------------------------------------------
# 2. Trim the copied sheet from 34 data rows down to the 26 the new challenge
#    needs, then restyle/re-word the tail rows that differ from a straight
#    truncation of Challenge 4.
# ---------------------------------------------------------------------------
$ws.Rows("27:34").Delete()

# -- Row 23: collapses the old "dashboards" sub-header row into a plain
#    (unlabeled) bullet-continuation row, and turns the merged note cell
#    into a single-row note.
$ws.Range("A23").ClearContents()
$ws.Range("F23").ClearContents()
$src.Range("F29").Copy()
$ws.Range("F23").PasteSpecial(-4122)

# -- Row 24: becomes a tall (34pt) bullet row using the "B36" note style.
$src.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").RowHeight = 34

# -- Row 25: plain bullet row (style already matches after the shift).

# -- Row 26: becomes the thick-bottomed final row of the section.
$src.Range("A34").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$src.Range("C34").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$src.Range("D34").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$src.Range("E34").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$src.Range("A34").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$src.Range("F28").Copy()
$ws.Range("F26").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Re-word the sheet for Challenge 5.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Challenge 5: Improve the environment"
$ws.Range("A16").Value = "Environment Improvements"
$ws.Range("C16").Formula = "=SUM(C17:C26)"
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = 29

$ws.Range("A17").Value = "How many improvements did the team implement that are:"
$ws.Range("F17").Value = "None: 0" + [char]10 + "1-2 Changes: 1" + [char]10 + "2-5 Changes: 2" + [char]10 + "5-9 Changes: 3" + [char]10 + "10+ Changes: 4"

$ws.Range("B23").Value = "All changes have been implemented via Azure Bicep/ARM Templates"
$ws.Range("B24").Value = "Application Insights has been added to the application and Grafana has been " + [char]10 + "updated to report the collected metrics"
$ws.Range("B25").Value = "Performance has been maintained or improved"
$ws.Range("B26").Value = "Accurately and effectively communicated the ACHIEVE process"

# ---------------------------------------------------------------------------
# 4. Column F width tweak that applies to every sheet in the workbook.
# ---------------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $col = $sheet.Range("F1").EntireColumn
    $col.ColumnWidth = 13.33203125
}
